# Daily attendance processing - 2026-01-28 15:15:52
#
# Normalizes the "Recorded By" column (G) on the "Session Analysis Results"
# sheet: whenever the automated "System" account is listed FIRST alongside a
# human reviewer's address, move "System" to the end of the list so the
# human reviewer is listed first, e.g.
#   "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#   "System, admin@admin.com"    -> "admin@admin.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count()

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    } elseif ($val -eq "System, admin@admin.com") {
        $cell.Value = "admin@admin.com, System"
    }
}
